# Fruta / hortaliza, semanal
# A new weekly record is inserted at row 406, pushing all the existing
# records (old rows 406..470) down by one row (to 407..471).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 406 (shifts rows 406..470 down to 407..471)
$ws.Rows.Item(406).Insert()

# Populate the newly inserted row 406 with the new weekly data point
$ws.Range("A406").Value = 6
$ws.Range("B406").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C406").Value = "Metropolitana"
$ws.Range("D406").Value = 45218
$ws.Range("E406").Value = 13
$ws.Range("F406").Value = 100112026
$ws.Range("G406").Value = "Haba"
$ws.Range("H406").Value = "Sin especificar"
$ws.Range("I406").Value = "Primera"
$ws.Range("J406").Value = 580
$ws.Range("K406").Value = 8000
$ws.Range("L406").Value = 9000
$ws.Range("M406").Value = 8552
$ws.Range("N406").Value = "$/saco 25 kilos"
$ws.Range("O406").Value = "Provincia de Melipilla"
$ws.Range("P406").Value = 342
$ws.Range("Q406").Value = 25
$ws.Range("R406").Value = "Hortaliza"
